$wb = $excel.ActiveWorkbook

# Rename "Sheet1" to "GSTInput-Cols"
$wsGst = $wb.Worksheets.Item("Sheet1")
$wsGst.Name = "GSTInput-Cols"

# On the "Data" sheet, update the selection to B1:D11 (active cell B1)
$wsData = $wb.Worksheets.Item("Data")
$wsData.Range("B1:D11").Select() | Out-Null

# Make "GSTInput-Cols" the active / selected tab, keeping its prior selection (F5)
$wsGst.Activate() | Out-Null
$wsGst.Range("F5").Select() | Out-Null
